$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - HYUNDAI ROTEM
$ws.Range("D2").Value = 181100
$ws.Range("E2").Value = 35.5
$ws.Range("F2").Value = 3.07
$ws.Range("N2").Value = 54.86376272656823

# Row 3 - HANWHA AEROSPACE
$ws.Range("D3").Value = 872500
$ws.Range("E3").Value = 32.3
$ws.Range("F3").Value = 2.41
$ws.Range("H3").Value = 60
$ws.Range("N3").Value = 54.86376272656823

# Row 4 - HANWHA SYSTEMS
$ws.Range("D4").Value = 46650
$ws.Range("E4").Value = 22.2
$ws.Range("F4").Value = 0.97
$ws.Range("N4").Value = 54.86376272656823

# Row 5 - KOREA AEROSPACE (only N5 changes)
$ws.Range("N5").Value = 54.86376272656823

# Row 6 - LIG Nex1
$ws.Range("D6").Value = 366500
$ws.Range("E6").Value = 26.7
$ws.Range("F6").Value = -4.31
$ws.Range("N6").Value = 54.86376272656823
